# Edit: add "valeur_portefeuille_buy_and_hold" column (I) and update
# "rendement_predit" (column B) with refreshed model values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update column B (rendement_predit) with the refreshed values ---
$bUpdates = @(
    @(3, 0.0006138469469600949),
    @(4, 0.06152842560246263),
    @(5, 0.05214948232460692),
    @(6, -0.04688672676649652),
    @(7, -0.01992836523031372),
    @(8, -0.03630076626129508),
    @(9, -0.006182425508914235),
    @(10, -0.006602968488977012),
    @(11, -0.02136562453104141),
    @(12, -0.01742764230501948),
    @(13, 0.02946957949093765),
    @(14, 0.02028816482165219),
    @(15, 0.004594878712980588),
    @(16, -0.007962266088021863),
    @(17, 0.002646617904192183),
    @(18, 0.01477172864263032),
    @(19, -0.005116722874085156),
    @(20, 0.008906113614404632),
    @(21, 0.02093297504097436),
    @(22, 0.00508627979317211),
    @(23, -0.00592431857048048),
    @(24, -0.01449278711531932),
    @(25, -0.02739007719521958),
    @(26, -0.004555390384473057),
    @(27, 0.002104718545581363),
    @(28, 0.005994587906439364),
    @(29, 0.001679560483630382),
    @(30, 0.01257061021838268),
    @(31, -0.007766238337204356),
    @(32, -0.00858756500723068),
    @(33, -0.00968254594673823),
    @(34, 0.0003146433812766247),
    @(35, 0.008503556081565833),
    @(36, 0.03792210973251109),
    @(37, 0.03280907582040626),
    @(38, -0.003937736148355242),
    @(39, -0.001792049111680072),
    @(40, 0.000408153491019192),
    @(41, -0.01982575703997647),
    @(42, -0.0191075814639472),
    @(43, -0.005841637504442332),
    @(44, 0.04350330945818648),
    @(45, -0.01385653001977616),
    @(46, -0.07985571145152903),
    @(47, -0.147797878767463),
    @(48, 0.01421689147663585),
    @(49, 0.006356375472220677),
    @(50, -0.009377149964747034),
    @(51, 0.005238368883798472),
    @(52, -0.02029794064702983),
    @(53, 0.004536062029579568),
    @(54, 0.005419553176341907),
    @(55, -0.04502589940245727),
    @(56, 0.006819967239302116),
    @(57, 0.02987242721610883),
    @(58, 0.03108308514259939),
    @(59, -0.001738685603724477),
    @(60, -0.03353070069447384),
    @(61, -0.00909949709778246),
    @(62, 0.03510229264497333),
    @(63, -0.008120757564272552),
    @(64, 0.01138770288412871),
    @(65, -0.01389573739630023),
    @(66, 0.004765423454299977),
    @(67, -0.01057497882278291),
    @(68, 0.01599020647839566),
    @(69, 0.002902061590788207),
    @(70, 0.01014394761165072),
    @(71, 0.01875378987535825),
    @(72, 0.01331820357692948),
    @(73, -0.01221526921111149),
    @(74, -0.04539433412887739),
    @(75, -0.02569790395934568),
    @(76, 0.01941483191194671),
    @(77, 0.007409066324184366),
    @(78, 0.007725413886689125),
    @(79, -0.005786295582733914),
    @(80, 0.02653668554586908),
    @(81, 0.003266878552935282),
    @(82, -0.02688463897347759),
    @(83, -0.00696056877955975),
    @(84, -0.0009128121504726749),
    @(85, 0.01401929550845082),
    @(86, 0.002203234996217773),
    @(87, -0.002011271541494608),
    @(88, -0.003383414973788135),
    @(89, 0.007426801281345874),
    @(90, 0.009610924328759651),
    @(91, 0.01115145561185216),
    @(92, 0.02314699186880453),
    @(93, 0.04101413864778891),
    @(94, 0.05199840802272959),
    @(95, 0.07899637171249019),
    @(96, 0.02685041389722898),
    @(97, -0.02179454713632722),
    @(98, 0.001687925836430182),
    @(99, 0.05118806656375874),
    @(100, 0.02612385304570886),
    @(101, -0.004649568397713466),
    @(102, 0.01544447260704729),
    @(103, 0.002546252367206847),
    @(104, 0.001479327953779475),
    @(106, 0.003790355670405532),
    @(107, 0.02532921722750103),
    @(108, 0.006713025280253859),
    @(109, -0.00001069802452846602),
    @(110, -0.02354498757245516),
    @(111, 0.007100425511689323),
    @(112, -0.006330658338711714),
    @(113, -0.04691550063389371),
    @(114, -0.0246605681318357),
    @(115, -0.01127476776701108),
    @(116, 0.02551659908697701),
    @(117, 0.07399982739437405),
    @(118, 0.003318976173316557),
    @(119, 0.03724807846133871),
    @(120, 0.02918870855696909),
    @(121, -0.006449298593647512),
    @(122, -0.02076956043908318),
    @(124, -0.02760987068615073),
    @(125, -0.007711237157845829),
    @(126, -0.01213762555422804),
    @(127, 0.01898610197877915),
    @(128, -0.001067065740514295),
    @(129, -0.03846943604605535),
    @(130, -0.01197838821331132),
    @(131, -0.008578513374683183),
    @(132, -0.01601527560768723),
    @(133, -0.05374448972491663),
    @(134, -0.02916788696566641),
    @(135, 0.1084222136179278),
    @(136, 0.07281846867300601),
    @(137, 0.008590415378815663),
    @(138, 0.02077642923662815),
    @(139, 0.06327023359835593),
    @(140, 0.03191158262469784),
    @(141, 0.02767680362445901),
    @(142, -0.01652256696269916),
    @(143, 0.01410901163084155)
)

foreach ($pair in $bUpdates) {
    $row = $pair[0]
    $val = $pair[1]
    $ws.Cells.Item($row, 2).Value = $val
}

# --- 2. Add new column I: "valeur_portefeuille_buy_and_hold" ---
# This is a simple buy-and-hold comparison portfolio: it starts from the same
# $10000 initial capital as the strategy (column G) but compounds using the
# raw observed market return (column F, rendement_observe) every day instead
# of the strategy's alpha-weighted return.
$ws.Cells.Item(1, 9).Value = "valeur_portefeuille_buy_and_hold"

# Match the formatting of the other header cells (bold, centered, bordered)
$ws.Cells.Item(1, 1).Copy()
$ws.Cells.Item(1, 9).PasteSpecial(-4122)

$lastRow = $ws.Cells.Item(1, 1).End(-4121).Row

$initialCapital = 10000
$previousValue = $initialCapital
for ($r = 2; $r -le $lastRow; $r++) {
    $observedReturn = $ws.Cells.Item($r, 6).Value()
    $currentValue = $previousValue * (1 + $observedReturn)
    $ws.Cells.Item($r, 9).Value = $currentValue
    $previousValue = $currentValue
}
